$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header date line: static date range -> templated placeholder ---
$ws.Range("A2").Value = "日期：{thisWeekStr}"

# --- "This week" schedule block (rows 10-14) ---
# Column A: weekday name placeholders (replacing literal 星期一..星期五)
# Column B: date placeholders (replacing literal date values), re-styled to
#           match column A's text style (numFmtId 49 / "@", same font/fill/border)

$ws.Range("A10").Value = "{thisWeek.Mon}`n"
$ws.Range("B10").Value = "{thisWeek.MonDate}`n"
$ws.Range("B10").NumberFormat = $ws.Range("A10").NumberFormat

$ws.Range("A11").Value = "{thisWeek.Tues}"
$ws.Range("B11").Value = "{thisWeek.TuesDate}"
$ws.Range("B11").NumberFormat = $ws.Range("A11").NumberFormat

$ws.Range("A12").Value = "{thisWeek.Wed}"
$ws.Range("B12").Value = "{thisWeek.WedDate}"
$ws.Range("B12").NumberFormat = $ws.Range("A12").NumberFormat
$ws.Range("C12").NumberFormat = "General"

$ws.Range("A13").Value = "{thisWeek.Thur}"
$ws.Range("B13").Value = "{thisWeek.ThurDate}"
$ws.Range("B13").NumberFormat = $ws.Range("A13").NumberFormat
$ws.Range("C13").NumberFormat = "General"

$ws.Range("A14").Value = "{thisWeek.Fri}"
$ws.Range("B14").Value = "{thisWeek.FriDate}"
$ws.Range("B14").NumberFormat = $ws.Range("A14").NumberFormat
$ws.Range("C14").NumberFormat = "General"

# --- "Next week" schedule block (rows 17-21) ---

$ws.Range("A17").Value = "{nextWeek.Mon}`n"
$ws.Range("B17").Value = "{nextWeek.MonDate}"
$ws.Range("B17").NumberFormat = $ws.Range("A17").NumberFormat

$ws.Range("A18").Value = "{nextWeek.Thes}`n"
$ws.Range("B18").Value = "{nextWeek.ThesDate}`n"
$ws.Range("B18").NumberFormat = $ws.Range("A18").NumberFormat

$ws.Range("A19").Value = "{nextWeek.Wed}`n"
$ws.Range("B19").Value = "{nextWeek.WedDate}`n"
$ws.Range("B19").NumberFormat = $ws.Range("A19").NumberFormat

$ws.Range("A20").Value = "{nextWeek.Thur}`n"
$ws.Range("B20").Value = "{nextWeek.ThurDate}`n"
$ws.Range("B20").NumberFormat = $ws.Range("A20").NumberFormat

$ws.Range("A21").Value = "{nextWeek.Fri}`n"
$ws.Range("B21").Value = "{nextWeek.FriDate}`n"
$ws.Range("B21").NumberFormat = $ws.Range("A21").NumberFormat
